$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C ("Variance" = B^2), shifting the old C ("Integral") to D
# and the old D ("Time") to E.
$ws.Columns.Item(3).Insert()

# Header + data for the new Variance column.
$ws.Range("C1").Value = "Variance"
$ws.Range("C2").Formula = "=B2^2"
$ws.Range("C3:C11").Formula = "=B3^2"

# Size column C to fit its new "Variance" header (best-effort AutoFit).
$ws.Columns.Item(3).ColumnWidth = 11.33

# Row 13 ("Avg") now also averages the new Variance column.
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Row 14 ("STD") no longer reports a standard deviation for the raw STD
# column (B) or the new Variance column (C) — only for the old Integral/Time
# columns (now D/E).
$ws.Range("B14:C14").ClearContents()

# New row 15: RMS computed as the square root of the averaged variance.
$ws.Range("A15").Value = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Match the author's final selection.
$null = $ws.Range("B15").Select()
